$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Part 1 dialogue update: translated flavor text / diaries / endings.
# Column C holds the new English translation next to the existing Japanese
# source text in column A. Shared strings must land in the workbook in the
# exact order Excel would append them (row 12 first, then 13, 18, 21, 5..11,
# 14..17, 19, 20) so write the cell values in that sequence.
# ---------------------------------------------------------------------------

$ws.Range("C12").Value2 = @"
\n<\n[1]>Right.\.\.
.\..\..\.\.
Let's go home.\.Meria.
"@

$ws.Range("C13").Value2 = "ED2　Succubus Hunter"

$ws.Range("C18").Value2 = @"
\n<\n[1]>...
"@

$ws.Range("C21").Value2 = @"
\n<\n[1]>Yeah.\.\.
\..\..\..\.\.
Let's go home.
"@

$ws.Range("C5").Value2 = @"
\n<Meria>Brother...
I'm so glad you're safe... Really...
Thank you for rescuing me...
"@

$ws.Range("C6").Value2 = @"
\n<Meria>If you hadn't come for me...
I never would have made it on my own.
"@

$ws.Range("C7").Value2 = @"
\n<Meria>Ummm...
I was completely broken down.
I was so terrified, and yet...
"@

$ws.Range("C8").Value2 = @"
\n<Meria>...
"@

$ws.Range("C9").Value2 = @"
\n<Meria>That's it!
"@

$ws.Range("C10").Value2 = @"
\n<Meria>A delicious feast!
I'll cook for you!
I have to stay strong! So...
"@

$ws.Range("C11").Value2 = @"
\n<Meria>.\..\..\.\.
Let's go.\.Big brother.\.\.
Back to our home.
"@

$ws.Range("C14").Value2 = @"
\n<Alsto>You were incredible... \n[1].
Thank you.
If it hadn't been for you, I'd be...
"@

$ws.Range("C15").Value2 = @"
\n<Alsto>You're a top notch hunter.
Much greater than me.
"@

$ws.Range("C16").Value2 = @"
\n<Alsto>I came here to help you, but the exact opposite happened.
I've lost all dignity as a big brother.
Or maybe I never had any in the first place!
"@

$ws.Range("C17").Value2 = @"
\n<Alsto>Anyway, weren't those 3 awfully strange?
I've never known succubi to work in groups.
And also...
"@

$ws.Range("C19").Value2 = @"
\n<Alsto>...
"@

$ws.Range("C20").Value2 = @"
\n<Alsto>No, it's probably nothing.\.\.
Let's go... Back to our home.
\n[1].
"@

# ---------------------------------------------------------------------------
# Formatting: wrap the multi-line translations, leave the short one-liners
# (the "..." beats and the ED2 title) un-wrapped but still on the new font,
# set the per-row heights to fit the wrapped text, and widen the columns.
# ---------------------------------------------------------------------------

$wrapRows = 5,6,7,10,11,12,14,15,16,17,20,21
foreach ($r in $wrapRows) {
    $ws.Range("C$r").WrapText = $true
}

$ws.Range("A6").WrapText = $true
$ws.Range("A16").WrapText = $true

$noWrapRows = 8,9,13,18,19
foreach ($r in $noWrapRows) {
    $ws.Range("C$r").Style = "Normal"
}

$rowHeights = @{5=45;6=45;7=45;10=45;11=45;12=45;14=45;15=30;16=75;17=45;20=45;21=45}
foreach ($r in $rowHeights.Keys) {
    $ws.Rows.Item($r).RowHeight = $rowHeights[$r]
}

$ws.Columns.Item(1).ColumnWidth = 48.333333333333336
$ws.Columns.Item(2).ColumnWidth = 51.666666666666664
$ws.Columns.Item(3).ColumnWidth = 62.833333333333336

$ws.Range("C25").Select()
